$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 497.64285
$ws.Range("I42").Value = 55.285713
$ws.Range("J42").Value = 940
$ws.Range("K42").Value = 165.857139
$ws.Range("L42").Value = 2820
$ws.Range("M42").Value = 64.14286099999998
$ws.Range("N42").Value = -3280

$ws.Range("H43").Value = 1786.9584
$ws.Range("I43").Value = 920.2
$ws.Range("J43").Value = 2015.0526
$ws.Range("K43").Value = 920.2
$ws.Range("L43").Value = 2015.0526
$ws.Range("M43").Value = -851.2
$ws.Range("N43").Value = -2153.0526

$ws.Range("H68").Value = 35295
$ws.Range("J68").Value = 35295
$ws.Range("L68").Value = 35295
$ws.Range("N68").Value = -36793

$ws.Range("H71").Value = 35295
$ws.Range("J71").Value = 35295
$ws.Range("L71").Value = 105885
$ws.Range("N71").Value = -113373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4679.7827
$ws.Range("I32").Value = 3451.0168
$ws.Range("K32").Value = 3451.0168
$ws.Range("M32").Value = -3164.0168

$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H103").Value = 25117.264
$ws.Range("J103").Value = 25117.264
$ws.Range("L103").Value = 25117.264
$ws.Range("N103").Value = -27461.264

$ws.Range("H110").Value = 5350
$ws.Range("I110").Value = 700
$ws.Range("J110").Value = 10000
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 10000
$ws.Range("M110").Value = 1345
$ws.Range("N110").Value = -14090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 150
$ws.Range("K22").Value = 150
$ws.Range("M22").Value = 23

$ws.Range("H54").Value = 2663.7144
$ws.Range("I54").Value = 2663.7144
$ws.Range("K54").Value = 2663.7144
$ws.Range("M54").Value = -2179.7144

$ws.Range("H80").Value = 602
$ws.Range("I80").Value = 636.1
$ws.Range("J80").Value = 571
$ws.Range("K80").Value = 636.1
$ws.Range("L80").Value = 571
$ws.Range("M80").Value = 361.9
$ws.Range("N80").Value = -2567

$ws.Range("H83").Value = 602
$ws.Range("I83").Value = 636.1
$ws.Range("J83").Value = 571
$ws.Range("K83").Value = 3180.5
$ws.Range("L83").Value = 2855
$ws.Range("M83").Value = 1811.5
$ws.Range("N83").Value = -12839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1054274.5
$ws.Range("I31").Value = 1563432.5
$ws.Range("J31").Value = 3109.5806
$ws.Range("K31").Value = 1563432.5
$ws.Range("L31").Value = 3109.5806
$ws.Range("M31").Value = -1563137.5
$ws.Range("N31").Value = -3699.5806

$ws.Range("H34").Value = 1054274.5
$ws.Range("I34").Value = 1563432.5
$ws.Range("J34").Value = 3109.5806
$ws.Range("K34").Value = 1563432.5
$ws.Range("L34").Value = 3109.5806
$ws.Range("M34").Value = -1563230.5
$ws.Range("N34").Value = -3513.5806

$ws.Range("H58").Value = 16669656
$ws.Range("I58").Value = 2085.0476
$ws.Range("J58").Value = 55560656
$ws.Range("K58").Value = 2085.0476
$ws.Range("L58").Value = 55560656
$ws.Range("M58").Value = -1882.0476
$ws.Range("N58").Value = -55561062

$ws.Range("H136").Value = 16669656
$ws.Range("I136").Value = 2085.0476
$ws.Range("J136").Value = 55560656
$ws.Range("K136").Value = 6255.1428
$ws.Range("L136").Value = 166681968
$ws.Range("M136").Value = -3705.1428
$ws.Range("N136").Value = -166687068

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2150.6155
$ws.Range("I64").Value = 1179
$ws.Range("J64").Value = 2983.4285
$ws.Range("K64").Value = 3537
$ws.Range("L64").Value = 8950.2855
$ws.Range("M64").Value = -3267
$ws.Range("N64").Value = -9490.2855

$ws.Range("H67").Value = 2150.6155
$ws.Range("I67").Value = 1179
$ws.Range("J67").Value = 2983.4285
$ws.Range("K67").Value = 3537
$ws.Range("L67").Value = 8950.2855
$ws.Range("M67").Value = -2601
$ws.Range("N67").Value = -10822.2855

$ws.Range("H68").Value = 2374.1667
$ws.Range("I68").Value = 713.9643
$ws.Range("J68").Value = 3826.8438
$ws.Range("K68").Value = 2141.8929
$ws.Range("L68").Value = 11480.5314
$ws.Range("M68").Value = -1330.8929
$ws.Range("N68").Value = -13102.5314

$ws.Range("H71").Value = 2374.1667
$ws.Range("I71").Value = 713.9643
$ws.Range("J71").Value = 3826.8438
$ws.Range("K71").Value = 6425.678699999999
$ws.Range("L71").Value = 34441.5942
$ws.Range("M71").Value = -2369.678699999999
$ws.Range("N71").Value = -42553.5942

$ws.Range("H131").Value = 1630.8
$ws.Range("I131").Value = 2490
$ws.Range("J131").Value = 1123.091
$ws.Range("K131").Value = 7470
$ws.Range("L131").Value = 3369.273
$ws.Range("M131").Value = -2430
$ws.Range("N131").Value = -13449.273

$ws.Range("H132").Value = 3239.4
$ws.Range("I132").Value = 2488.2222
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 22393.9998
$ws.Range("L132").Value = 90000
$ws.Range("M132").Value = -19863.9998
$ws.Range("N132").Value = -95060

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 80042
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 80042
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 80042
$ws.Range("N26").Value = -80602
$ws.Range("M26").ClearContents()

$ws.Range("H50").Value = 80042
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 80042
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 80042
$ws.Range("N50").Value = -81038
$ws.Range("M50").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 166669360
$ws.Range("I22").Value = 500000220
$ws.Range("K22").Value = 500000220
$ws.Range("M22").Value = -499999925

$ws.Range("H27").Value = 166669360
$ws.Range("I27").Value = 500000220
$ws.Range("K27").Value = 500000220
$ws.Range("M27").Value = -500000113

$ws.Range("H100").Value = 1966.5714
$ws.Range("I100").Value = 827.4286
$ws.Range("J100").Value = 3105.7144
$ws.Range("K100").Value = 827.4286
$ws.Range("L100").Value = 3105.7144
$ws.Range("M100").Value = -286.4286
$ws.Range("N100").Value = -4187.7144

$ws.Range("H122").Value = 3058.7942
$ws.Range("I122").Value = 2551.6897
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 7655.0691
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5205.0691
$ws.Range("N122").Value = -22900

$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
